# Automatische test-sync: 2025-06-27 22:47:50
#
# Appends a new "Productinformatie" test-mail entry (mail #3) to the Logs
# sheet, rolls the matching category tally into the Dashboard sheet, and
# extends the conditional formatting / chart ranges so they keep covering
# the newly added row.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# ---------------------------------------------------------------------
# 1. Append the new row to the "Logs" sheet (row 12)
# ---------------------------------------------------------------------
$newLogRow = 12

$logs.Range("A" + $newLogRow).Value = "Wat zijn de verzendkosten?"
$logs.Range("B" + $newLogRow).Value = "mailmind.test@zohomail.eu"
$logs.Range("C" + $newLogRow).Value = "Testmail #3: Wat zijn de verzendkosten?"
$logs.Range("D" + $newLogRow).Value = "Productinformatie"
$logs.Range("E" + $newLogRow).Value = "Beste afzender,`nDank u voor uw interesse in onze producten/diensten. Om u nauwkeurige informatie te verstrekken over de verzendkosten, hebben we meer details nodig zoals het product/dienst waar u naar informeert en het afleveradres. Zou u ons kunnen voorzien van deze informatie zodat we u een precieze schatting van de verzendkosten kunnen geven?`nMet vriendelijke groet,`n[Naam]`nE-mailassistent"
$logs.Range("F" + $newLogRow).Value = "2025-06-27 22:46:58"
$logs.Range("G" + $newLogRow).Value = "Ja"
$logs.Range("H" + $newLogRow).Value = "Nee"
$logs.Range("I" + $newLogRow).Value = "Ja"

# The description field contains embedded line breaks, which makes the
# emulated engine mark the row with an explicit (auto-fit) height. Real
# Excel leaves the row at the default/standard height here, so re-fit it
# back to standard to keep the row definition clean.
$logs.Rows.Item($newLogRow).AutoFit()

# ---------------------------------------------------------------------
# 2. Extend the conditional formatting ranges on "Logs" from row 11 to
#    row 12 for every column that carries formatting (D, G, H, I), while
#    preserving rule order / priorities / linked styles.
# ---------------------------------------------------------------------
$oldLastRow = 11
$newLastRow = 12

foreach ($column in @("D", "G", "H", "I")) {
    $oldRange = $logs.Range($column + "2:" + $column + $oldLastRow)
    $newRange = $logs.Range($column + "2:" + $column + $newLastRow)
    $rules = $oldRange.FormatConditions

    for ($i = 1; $i -le $rules.Count; $i++) {
        $rules.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------------
# 3. Append the matching tally row to the "Dashboard" sheet (row 6)
# ---------------------------------------------------------------------
$newDashboardRow = 6

$dashboard.Range("A" + $newDashboardRow).Value = "Productinformatie"
$dashboard.Range("B" + $newDashboardRow).Value = 1

# ---------------------------------------------------------------------
# 4. Extend the chart's category/value series references to include the
#    new Dashboard row.
# ---------------------------------------------------------------------
$chartObject = $dashboard.ChartObjects().Item(1)
$series = $chartObject.Chart.SeriesCollection(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$" + $newDashboardRow + ",'Dashboard'!`$B`$2:`$B`$" + $newDashboardRow + ",1)"
